{"js": "// Update the two-digit \u00f7 one-digit division answer table cells to the\n// new set of problems/answers (see commit \"Update master to output generated at c8c62b6\").\nconst replacements = [\n  [\"79\u00f79=8, 7\", \"70\u00f74=17, 2\"],\n  [\"69\u00f77=9, 6\", \"71\u00f77=10, 1\"],\n  [\"45\u00f78=5, 5\", \"52\u00f74=13, 0\"],\n  [\"93\u00f74=23, 1\", \"97\u00f72=48, 1\"],\n  [\"27\u00f75=5, 2\", \"41\u00f72=20, 1\"],\n  [\"66\u00f73=22, 0\", \"74\u00f72=37, 0\"],\n  [\"51\u00f76=8, 3\", \"31\u00f72=15, 1\"],\n  [\"45\u00f73=15, 0\", \"42\u00f73=14, 0\"],\n  [\"45\u00f75=9, 0\", \"73\u00f79=8, 1\"],\n  [\"25\u00f73=8, 1\", \"71\u00f79=7, 8\"],\n  [\"13\u00f77=1, 6\", \"48\u00f78=6, 0\"],\n  [\"23\u00f72=11, 1\", \"56\u00f78=7, 0\"],\n  [\"95\u00f72=47, 1\", \"57\u00f74=14, 1\"],\n  [\"49\u00f75=9, 4\", \"24\u00f73=8, 0\"],\n  [\"20\u00f73=6, 2\", \"89\u00f77=12, 5\"],\n  [\"10\u00f76=1, 4\", \"80\u00f76=13, 2\"],\n  [\"66\u00f75=13, 1\", \"27\u00f72=13, 1\"],\n  [\"52\u00f73=17, 1\", \"69\u00f72=34, 1\"],\n  [\"26\u00f72=13, 0\", \"26\u00f79=2, 8\"],\n  [\"44\u00f79=4, 8\", \"75\u00f78=9, 3\"],\n  [\"34\u00f77=4, 6\", \"58\u00f75=11, 3\"],\n  [\"71\u00f78=8, 7\", \"10\u00f78=1, 2\"],\n  [\"88\u00f77=12, 4\", \"93\u00f73=31, 0\"],\n  [\"62\u00f77=8, 6\", \"34\u00f73=11, 1\"],\n  [\"65\u00f79=7, 2\", \"26\u00f75=5, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit \u00f7 one-digit division answer table cells to the\n# new set of problems/answers (see commit \"Update master to output generated at c8c62b6\").\n$d = $word.ActiveDocument\n\n$null = $d.Content.Find.Execute(\"79\u00f79=8, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"70\u00f74=17, 2\", 2)\n$null = $d.Content.Find.Execute(\"69\u00f77=9, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"71\u00f77=10, 1\", 2)\n$null = $d.Content.Find.Execute(\"45\u00f78=5, 5\", $false, $false, $false, $false, $false, $true, 1, $false, \"52\u00f74=13, 0\", 2)\n$null = $d.Content.Find.Execute(\"93\u00f74=23, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"97\u00f72=48, 1\", 2)\n$null = $d.Content.Find.Execute(\"27\u00f75=5, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"41\u00f72=20, 1\", 2)\n$null = $d.Content.Find.Execute(\"66\u00f73=22, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"74\u00f72=37, 0\", 2)\n$null = $d.Content.Find.Execute(\"51\u00f76=8, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"31\u00f72=15, 1\", 2)\n$null = $d.Content.Find.Execute(\"45\u00f73=15, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"42\u00f73=14, 0\", 2)\n$null = $d.Content.Find.Execute(\"45\u00f75=9, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"73\u00f79=8, 1\", 2)\n$null = $d.Content.Find.Execute(\"25\u00f73=8, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"71\u00f79=7, 8\", 2)\n$null = $d.Content.Find.Execute(\"13\u00f77=1, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"48\u00f78=6, 0\", 2)\n$null = $d.Content.Find.Execute(\"23\u00f72=11, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00f78=7, 0\", 2)\n$null = $d.Content.Find.Execute(\"95\u00f72=47, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"57\u00f74=14, 1\", 2)\n$null = $d.Content.Find.Execute(\"49\u00f75=9, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"24\u00f73=8, 0\", 2)\n$null = $d.Content.Find.Execute(\"20\u00f73=6, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"89\u00f77=12, 5\", 2)\n$null = $d.Content.Find.Execute(\"10\u00f76=1, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"80\u00f76=13, 2\", 2)\n$null = $d.Content.Find.Execute(\"66\u00f75=13, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"27\u00f72=13, 1\", 2)\n$null = $d.Content.Find.Execute(\"52\u00f73=17, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"69\u00f72=34, 1\", 2)\n$null = $d.Content.Find.Execute(\"26\u00f72=13, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"26\u00f79=2, 8\", 2)\n$null = $d.Content.Find.Execute(\"44\u00f79=4, 8\", $false, $false, $false, $false, $false, $true, 1, $false, \"75\u00f78=9, 3\", 2)\n$null = $d.Content.Find.Execute(\"34\u00f77=4, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"58\u00f75=11, 3\", 2)\n$null = $d.Content.Find.Execute(\"71\u00f78=8, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"10\u00f78=1, 2\", 2)\n$null = $d.Content.Find.Execute(\"88\u00f77=12, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"93\u00f73=31, 0\", 2)\n$null = $d.Content.Find.Execute(\"62\u00f77=8, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"34\u00f73=11, 1\", 2)\n$null = $d.Content.Find.Execute(\"65\u00f79=7, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"26\u00f75=5, 1\", 2)\n"}
